# Auto-generated edit script: rebuild the Data/Tags colnames table (A1:D25)
# matching the updated "Bibliometric_workshop" column-name reference sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous 30-row table; the refreshed table only has 25 rows
# (JI/SC/UT/WC/Z9 tag rows removed, PY promoted, PF description reworded).
$ws.Range("A1:D30").ClearContents()

$ws.Cells.Item(1, 1).Value = 'Index'
$ws.Cells.Item(1, 2).Value = 'Var name'
$ws.Cells.Item(1, 3).Value = 'Description'
$ws.Cells.Item(1, 4).Value = 'Var Type'

$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 'AU'
$ws.Cells.Item(2, 3).Value = 'Authors'
$ws.Cells.Item(2, 4).Value = 'Text'

$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 'TI'
$ws.Cells.Item(3, 3).Value = 'Document Title'
$ws.Cells.Item(3, 4).Value = 'Text'

$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 'SO'
$ws.Cells.Item(4, 3).Value = 'Publication Name'
$ws.Cells.Item(4, 4).Value = 'Text'

$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 'PY'
$ws.Cells.Item(5, 3).Value = 'Year Published'
$ws.Cells.Item(5, 4).Value = 'Numeric'

$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 'DT'
$ws.Cells.Item(6, 3).Value = 'Document Type'
$ws.Cells.Item(6, 4).Value = 'Text'

$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 'DE'
$ws.Cells.Item(7, 3).Value = 'Author Keywords'
$ws.Cells.Item(7, 4).Value = 'Text'

$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 'ID'
$ws.Cells.Item(8, 3).Value = 'Keywords Plus®'
$ws.Cells.Item(8, 4).Value = 'Text'

$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 'AB'
$ws.Cells.Item(9, 3).Value = 'Abstract'
$ws.Cells.Item(9, 4).Value = 'Text'

$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 'C1'
$ws.Cells.Item(10, 3).Value = 'Author Address'
$ws.Cells.Item(10, 4).Value = 'Text'

$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 'RP'
$ws.Cells.Item(11, 3).Value = 'Reprint Address'
$ws.Cells.Item(11, 4).Value = 'Text'

$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = 'CR'
$ws.Cells.Item(12, 3).Value = 'Cited References'
$ws.Cells.Item(12, 4).Value = 'Text'

$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = 'TC'
$ws.Cells.Item(13, 3).Value = 'Web of Science Core Collection Times Cited Count'
$ws.Cells.Item(13, 4).Value = 'Numeric'

$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = 'DB'
$ws.Cells.Item(14, 3).Value = 'Data Base'
$ws.Cells.Item(14, 4).Value = 'Text'

$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = 'AU_UN'
$ws.Cells.Item(15, 3).Value = 'Authors affiliations'
$ws.Cells.Item(15, 4).Value = 'Text'

$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = 'AU1_UN'
$ws.Cells.Item(16, 3).Value = 'Corresponding Author affiliation'
$ws.Cells.Item(16, 4).Value = 'Text'

$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = 'AU_UN_NR'
$ws.Cells.Item(17, 4).Value = 'Text'

$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = 'SR_FULL'
$ws.Cells.Item(18, 3).Value = 'Author year journal '
$ws.Cells.Item(18, 4).Value = 'Text'

$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = 'SR'
$ws.Cells.Item(19, 3).Value = 'Author year journal full'
$ws.Cells.Item(19, 4).Value = 'Text'

$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = 'AU_CO'
$ws.Cells.Item(20, 3).Value = 'Authors country'
$ws.Cells.Item(20, 4).Value = 'Text'

$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = 'DI'
$ws.Cells.Item(21, 3).Value = 'Digital Object Identifier (DOI)'
$ws.Cells.Item(21, 4).Value = 'Text'

$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = 'BN'
$ws.Cells.Item(22, 3).Value = 'International Standard Book Number (ISBN)'
$ws.Cells.Item(22, 4).Value = 'Text'

$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = 'PU'
$ws.Cells.Item(23, 3).Value = 'Publisher'
$ws.Cells.Item(23, 4).Value = 'Text'

$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = 'FU'
$ws.Cells.Item(24, 3).Value = 'Funding Agency and Grant Number'
$ws.Cells.Item(24, 4).Value = 'Text'

$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = 'PF'
$ws.Cells.Item(25, 3).Value = 'Platform (Web of Science or Scopus)'
$ws.Cells.Item(25, 4).Value = 'Text'

# Match the author's final selection/cursor position on the sheet
$ws.Range("D26").Select()
